# MOS-23045: Update Master Data as per 22 April Changes
# Appends 10 new "Postal Code" hierarchy rows (eng/fra/ara) under the
# BNMR (Ben Mansour) parent location to the master-location sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-location")

# Columns: A=code, B=name, C=hierarchy_level, D=hierarchy_level_name,
#          E=parent_loc_code, F=lang_code, G=is_active, H=cr_by, I=cr_dtimes
$newRows = @(
    @(10110, 10110, 5, "Postal Code", "BNMR", "eng"),
    @(10111, 10111, 5, "Postal Code", "BNMR", "eng"),
    @(10113, 10113, 5, "Postal Code", "BNMR", "eng"),
    @(10114, 10114, 5, "Postal Code", "BNMR", "eng"),
    @(10111, 10111, 5, "code postal", "BNMR", "fra"),
    @(10110, 10110, 5, "code postal", "BNMR", "fra"),
    @(10113, 10113, 5, "code postal", "BNMR", "fra"),
    @(10114, 10114, 5, "code postal", "BNMR", "fra"),
    @(10111, 10111, 5, "الرمز البريدي", "BNMR", "ara"),
    @(10110, 10110, 5, "الرمز البريدي", "BNMR", "ara")
)

$startRow = 110
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $true
    $ws.Cells.Item($r, 8).Value = "superadmin"
    $ws.Cells.Item($r, 9).Value = "now()"
}

# Reset the lingering selection (previously A110:XFD112) to a single cell.
$ws.Range("A1").Select()
